$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: (A-AVERAGE(A2:A11))*(B-AVERAGE(B2:B11))
$ws.Range("C2").Formula = "=(A2-AVERAGE(`$A`$2:`$A`$11))*(B2-AVERAGE(`$B`$2:`$B`$11))"
$ws.Range("C3:C11").Formula = "=(A3-AVERAGE(`$A`$2:`$A`$11))*(B3-AVERAGE(`$B`$2:`$B`$11))"

# Column E: covariance functions
$ws.Range("E3").Formula = "=_xlfn.COVARIANCE.P(A2:A11, B2:B11)"
$ws.Range("E4").Formula = "=_xlfn.COVARIANCE.S(A2:A11, B2:B11)"

# C12: average of C2:C11
$ws.Range("C12").Formula = "=AVERAGE(C2:C11)"

# Selection changes to D2
$ws.Range("D2").Select()

# Column D width (stored width of 9 chars; COM ColumnWidth has a constant ~5/6 offset
# vs. the stored OOXML column width in this environment)
$ws.Columns.Item(4).ColumnWidth = 8.166666666666666
